$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 450-532: Fecha(D), Calidad(L), Volumen(M), PrecioMin(N), PrecioMax(O), PrecioProm(P), Origen(R), PrecioKg(S)
$data = @(
    @(44995, "Especial", 200, 14000, 15000, 14500, "Provincia de Melipilla", 806),
    @(44995, "Extra (doble especial)", 180, 16000, 16000, 16000, "Provincia de Melipilla", 889),
    @(44995, "Primera", 250, 12000, 13000, 12500, "Provincia de Melipilla", 694),
    @(44995, "Segunda", 150, 8000, 10000, 9000, "Provincia de Melipilla", 500),
    @(44995, "Tercera", 200, 7000, 7000, 7000, "Provincia de Melipilla", 389),
    @(44798, "Especial", 50, 25000, 25000, 25000, "Provincia de Melipilla", 1389),
    @(44798, "Extra (doble especial)", 70, 30000, 30000, 30000, "Provincia de Melipilla", 1667),
    @(44798, "Primera", 40, 20000, 20000, 20000, "Provincia de Melipilla", 1111),
    @(44274, "Especial", 200, 14000, 15000, 14500, "Provincia de Melipilla", 806),
    @(44274, "Extra (doble especial)", 220, 16000, 17000, 16500, "Provincia de Melipilla", 917),
    @(44274, "Primera", 240, 12000, 13000, 12500, "Provincia de Melipilla", 694),
    @(44286, "Especial", 150, 15000, 15000, 15000, "Provincia de Melipilla", 833),
    @(44286, "Extra (doble especial)", 120, 17000, 17000, 17000, "Provincia de Melipilla", 944),
    @(44286, "Primera", 145, 13000, 13000, 13000, "Provincia de Melipilla", 722),
    @(44286, "Segunda", 100, 10000, 10000, 10000, "Provincia de Melipilla", 556),
    @(44721, "Cuarta", 40, 10000, 10000, 10000, "Provincia de Melipilla", 556),
    @(44721, "Especial", 80, 23000, 23000, 23000, "Provincia de Melipilla", 1278),
    @(44721, "Primera", 75, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44655, "Especial", 100, 12000, 12000, 12000, "Provincia de Limarí", 667),
    @(44655, "Primera", 170, 9000, 9000, 9000, "Provincia de Limarí", 500),
    @(44655, "Segunda", 150, 7000, 7000, 7000, "Provincia de Limarí", 389),
    @(44974, "Especial", 120, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44974, "Extra (doble especial)", 70, 20000, 20000, 20000, "Provincia de Melipilla", 1111),
    @(44974, "Primera", 260, 16000, 16000, 16000, "Provincia de Melipilla", 889),
    @(44974, "Segunda", 115, 13000, 13000, 13000, "Provincia de Melipilla", 722),
    @(44963, "Primera", 275, 22000, 22000, 22000, "Región Metropolitana", 1222),
    @(44963, "Segunda", 275, 18000, 18000, 18000, "Región Metropolitana", 1000),
    @(44963, "Tercera", 275, 15000, 15000, 15000, "Región Metropolitana", 833),
    @(44966, "Especial", 200, 20000, 20000, 20000, "Región Metropolitana", 1111),
    @(44966, "Primera", 200, 17000, 17000, 17000, "Región Metropolitana", 944),
    @(44966, "Segunda", 200, 13000, 13000, 13000, "Región Metropolitana", 722),
    @(44987, "Primera", 300, 13000, 14000, 13500, "Región Metropolitana", 750),
    @(44987, "Segunda", 265, 9000, 9000, 9000, "Región Metropolitana", 500),
    @(44672, "Especial", 150, 16000, 16000, 16000, "Provincia de Melipilla", 889),
    @(44672, "Extra (doble especial)", 120, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44672, "Primera", 100, 14000, 14000, 14000, "Provincia de Melipilla", 778),
    @(44672, "Segunda", 150, 12000, 12000, 12000, "Provincia de Melipilla", 667),
    @(44383, "Primera", 200, 14000, 14000, 14000, "Provincia de Santiago", 778),
    @(44383, "Segunda", 200, 12000, 12000, 12000, "Provincia de Santiago", 667),
    @(44383, "Tercera", 170, 8000, 8000, 8000, "Provincia de Santiago", 444),
    @(44608, "Especial", 135, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44608, "Primera", 150, 15000, 15000, 15000, "Provincia de Melipilla", 833),
    @(44608, "Segunda", 60, 12000, 12000, 12000, "Provincia de Melipilla", 667),
    @(44511, "Especial", 150, 33000, 33000, 33000, "Provincia de Melipilla", 1833),
    @(44511, "Primera", 20, 28000, 28000, 28000, "Provincia de Melipilla", 1556),
    @(44511, "Segunda", 25, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44720, "Especial", 80, 22000, 22000, 22000, "Provincia de Melipilla", 1222),
    @(44720, "Primera", 100, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44665, "Especial", 125, 16000, 16000, 16000, "Provincia de Melipilla", 889),
    @(44665, "Extra (doble especial)", 100, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44665, "Primera", 150, 14000, 14000, 14000, "Provincia de Melipilla", 778),
    @(44665, "Segunda", 150, 12000, 12000, 12000, "Provincia de Melipilla", 667),
    @(44270, "Especial", 150, 14000, 14000, 14000, "Provincia de Melipilla", 778),
    @(44270, "Extra (doble especial)", 120, 17000, 17000, 17000, "Provincia de Melipilla", 944),
    @(44270, "Primera", 180, 12000, 12000, 12000, "Provincia de Melipilla", 667),
    @(44270, "Segunda", 100, 9000, 9000, 9000, "Provincia de Melipilla", 500),
    @(44273, "Especial", 150, 15000, 15000, 15000, "Provincia de Melipilla", 833),
    @(44273, "Extra (doble especial)", 100, 17000, 17000, 17000, "Provincia de Melipilla", 944),
    @(44273, "Primera", 180, 12000, 12000, 12000, "Provincia de Melipilla", 667),
    @(44273, "Segunda", 120, 10000, 10000, 10000, "Provincia de Melipilla", 556),
    @(44273, "Tercera", 125, 8000, 8000, 8000, "Provincia de Melipilla", 444),
    @(44257, "Especial", 240, 17000, 18000, 17500, "Provincia de Melipilla", 972),
    @(44257, "Extra (doble especial)", 80, 20000, 20000, 20000, "Provincia de Melipilla", 1111),
    @(44257, "Primera", 275, 14000, 14000, 14000, "Provincia de Melipilla", 778),
    @(44257, "Segunda", 230, 10000, 12000, 10696, "Provincia de Melipilla", 594),
    @(44670, "Especial", 50, 16000, 16000, 16000, "Provincia de Melipilla", 889),
    @(44670, "Extra (doble especial)", 100, 18000, 18000, 18000, "Provincia de Melipilla", 1000),
    @(44670, "Primera", 30, 14000, 14000, 14000, "Provincia de Melipilla", 778),
    @(44670, "Segunda", 20, 12000, 12000, 12000, "Provincia de Melipilla", 667),
    @(44278, "Especial", 280, 14000, 15000, 14286, "Provincia de Melipilla", 794),
    @(44278, "Primera", 210, 13000, 13000, 13000, "Provincia de Melipilla", 722),
    @(44278, "Segunda", 80, 10000, 11000, 10625, "Provincia de Melipilla", 590),
    @(44278, "Tercera", 50, 9000, 9000, 9000, "Provincia de Melipilla", 500),
    @(44432, "Especial", 50, 35000, 35000, 35000, "Provincia de Melipilla", 1944),
    @(44432, "Extra (doble especial)", 70, 38000, 38000, 38000, "Provincia de Melipilla", 2111),
    @(44432, "Primera", 10, 20000, 20000, 20000, "Provincia de Melipilla", 1111),
    @(44432, "Segunda", 5, 15000, 15000, 15000, "Provincia de Melipilla", 833),
    @(44399, "Primera", 15, 20000, 20000, 20000, "Provincia de Melipilla", 1111),
    @(44399, "Segunda", 15, 15000, 15000, 15000, "Provincia de Melipilla", 833),
    @(44399, "Tercera", 50, 10000, 10000, 10000, "Provincia de Melipilla", 556),
    @(44658, "Especial", 170, 11000, 11000, 11000, "Región Metropolitana", 611),
    @(44658, "Primera", 250, 8000, 8000, 8000, "Región Metropolitana", 444),
    @(44658, "Segunda", 230, 5000, 5000, 5000, "Región Metropolitana", 278)
)

$startRow = 450
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Rows beyond the previous last row (527) are brand new rows; fill constant columns too
    if ($r -gt 527) {
        $ws.Cells.Item($r, 1).Value = 6
        $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
        $ws.Cells.Item($r, 3).Value = "Metropolitana"
        $ws.Cells.Item($r, 5).Value = 13
        $ws.Cells.Item($r, 6).Value = "Fruta"
        $ws.Cells.Item($r, 7).Value = 100107
        $ws.Cells.Item($r, 8).Value = "Otros"
        $ws.Cells.Item($r, 9).Value = 100107011
        $ws.Cells.Item($r, 10).Value = "Tuna"
        $ws.Cells.Item($r, 11).Value = "Sin especificar"
        $ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
        $ws.Cells.Item($r, 20).Value = 18
        $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(527, 4).NumberFormat
    }

    $ws.Cells.Item($r, 4).Value = $row[0]
    $ws.Cells.Item($r, 12).Value = $row[1]
    $ws.Cells.Item($r, 13).Value = $row[2]
    $ws.Cells.Item($r, 14).Value = $row[3]
    $ws.Cells.Item($r, 15).Value = $row[4]
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 18).Value = $row[6]
    $ws.Cells.Item($r, 19).Value = $row[7]
}

Write-Output "Update complete"
